# Apply the updated crypto price/volume snapshot values.
# Cells whose new text looks like a plain number (e.g. "309.64") are written
# with a leading apostrophe so Excel keeps them as text, matching the source data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '44.651.88'
$ws.Range('E2').Value = '  +1.48%  '
$ws.Range('D3').Value = '2.253.46'
$ws.Range('E3').Value = '  +0.81%  '
$ws.Range('E4').Value = '  +0.91%  '
$ws.Range('D5').Value = '''309.64'  # force text
$ws.Range('E5').Value = '  +1.52%  '
$ws.Range('D6').Value = '''95.60'  # force text
$ws.Range('E6').Value = '  +0.57%  '
$ws.Range('E7').Value = '  +1.40%  '
$ws.Range('E8').Value = '  +0.23%  '
$ws.Range('D9').Value = '''0.530'  # force text
$ws.Range('E9').Value = '  +1.55%  '
$ws.Range('D10').Value = '''35.26'  # force text
$ws.Range('E10').Value = '  +2.21%  '
$ws.Range('E11').Value = '  +0.93%  '
$ws.Range('D12').Value = '''7.33'  # force text
$ws.Range('E12').Value = '  +2.17%  '
$ws.Range('E13').Value = '  +0.74%  '
$ws.Range('E14').Value = '  +3.65%  '
$ws.Range('D15').Value = '2.260.64'
$ws.Range('E15').Value = '  +1.08%  '
$ws.Range('D16').Value = '''13.73'  # force text
$ws.Range('E16').Value = '  +2.02%  '
$ws.Range('D17').Value = '44.334.75'
$ws.Range('E17').Value = '  +1.07%  '
$ws.Range('D18').Value = '0.0₃0967'
$ws.Range('E18').Value = '  +1.24%  '
$ws.Range('B19').Value = 'InternetComputer(DFINITY)'
$ws.Range('C19').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D19').Value = '''12.27'  # force text
$ws.Range('E19').Value = '  +1.11%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').Value = '''6.43'  # force text
$ws.Range('E20').Value = '  +4.40%  '
$ws.Range('D21').Value = '''66.02'  # force text
$ws.Range('E21').Value = '  +2.10%  '
$ws.Range('D22').Value = '''240.24'  # force text
$ws.Range('E22').Value = '  +1.24%  '
$ws.Range('E23').Value = '  +3.41%  '
$ws.Range('D24').Value = '''2.01'  # force text
$ws.Range('E24').Value = '  +3.95%  '
$ws.Range('E25').Value = '  -0.38%  '
$ws.Range('E26').Value = '  +5.10%  '
$ws.Range('D27').Value = '''9.89'  # force text
$ws.Range('E27').Value = '  +0.63%  '
$ws.Range('D28').Value = '''37.86'  # force text
$ws.Range('E28').Value = '  +3.56%  '
$ws.Range('D29').Value = '''6.06'  # force text
$ws.Range('E29').Value = '  +3.47%  '
$ws.Range('E30').Value = '  +0.76%  '
$ws.Range('D31').Value = '''152.94'  # force text
$ws.Range('E31').Value = '  -0.50%  '
$ws.Range('D32').Value = '''0.0811'  # force text
$ws.Range('E32').Value = '  +0.86%  '
$ws.Range('E33').Value = '  +1.01%  '
$ws.Range('D34').Value = '''3.18'  # force text
$ws.Range('E34').Value = '  -1.47%  '
$ws.Range('D35').Value = '''0.111'  # force text
$ws.Range('E35').Value = '  +0.65%  '
$ws.Range('E36').Value = '  +2.07%  '
$ws.Range('D37').Value = '''1.82'  # force text
$ws.Range('E37').Value = '  +3.48%  '
$ws.Range('E38').Value = '  +4.19%  '
$ws.Range('D39').Value = '''3.84'  # force text
$ws.Range('E39').Value = '  +1.92%  '
$ws.Range('D40').Value = '''14.49'  # force text
$ws.Range('E40').Value = '  -2.19%  '
$ws.Range('E41').Value = '  +1.45%  '
$ws.Range('E42').Value = '  +0.34%  '
$ws.Range('D43').Value = '1.756.59'
$ws.Range('E43').Value = '  +1.51%  '
$ws.Range('E44').Value = '  +5.45%  '
$ws.Range('D45').Value = '''81.50'  # force text
$ws.Range('E45').Value = '  -3.25%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').Value = '''100.31'  # force text
$ws.Range('E46').Value = '  +0.73%  '
$ws.Range('B47').Value = 'ordi'
$ws.Range('C47').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D47').Value = '''71.30'  # force text
$ws.Range('E47').Value = '  +4.22%  '
$ws.Range('B48').Value = 'THORChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D48').Value = '''4.91'  # force text
$ws.Range('E48').Value = '  +0.86%  '
$ws.Range('B49').Value = 'MultiversX'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D49').Value = '''55.82'  # force text
$ws.Range('E49').Value = '  +3.55%  '
$ws.Range('B50').Value = 'FraxShare'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D50').Value = '''8.22'  # force text
$ws.Range('E50').Value = '  +2.40%  '
$ws.Range('E51').Value = '  +5.42%  '
